# Auto-generated edit script for cryptos.xlsx update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 38/39: coin name + link swap (Fetch.AI <-> Bittensor) ---
$ws.Range("B38").Value2 = 'Bittensor'
$ws.Range("C38").Value2 = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("B39").Value2 = 'Fetch.AI'
$ws.Range("C39").Value2 = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'

# --- Price (column D) updates; force text storage, then strip the temporary
#     number-format style so cells keep their original (unstyled) appearance ---
$priceUpdates = @{
    "D2" = '70.793.57'
    "D3" = '3.630.02'
    "D4" = '0.998'
    "D5" = '605.11'
    "D6" = '199.51'
    "D7" = '0.627'
    "D8" = '0.999'
    "D9" = '0.219'
    "D10" = '0.646'
    "D11" = '53.78'
    "D12" = '0.0000306'
    "D13" = '9.55'
    "D14" = '4.193.26'
    "D15" = '612.84'
    "D16" = '13.02'
    "D17" = '70.837.13'
    "D18" = '3.631.52'
    "D19" = '19.06'
    "D21" = '0.999'
    "D22" = '18.29'
    "D23" = '5.36'
    "D24" = '103.60'
    "D25" = '4.63'
    "D26" = '3.00'
    "D27" = '10.60'
    "D28" = '9.75'
    "D29" = '33.72'
    "D30" = '4.74'
    "D31" = '7.20'
    "D32" = '12.26'
    "D34" = '63.36'
    "D36" = '3.986.10'
    "D38" = '516.20'
    "D39" = '3.04'
    "D40" = '0.390'
    "D41" = '36.69'
    "D45" = '3.50'
    "D47" = '0.141'
    "D48" = '8.60'
    "D50" = '0.000250'
}
foreach ($addr in $priceUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value2 = $priceUpdates[$addr]
    $cell.ClearFormats()
}

# --- Volume(1h) (column E) updates; these are already non-numeric text
#     (leading/trailing spaces + '%'), so a plain assignment keeps them as text ---
$volumeUpdates = @{
    "E2" = '  +1.72%  '
    "E3" = '  +3.91%  '
    "E4" = '  -0.09%  '
    "E5" = '  +0.36%  '
    "E6" = '  +3.02%  '
    "E7" = '  +0.33%  '
    "E8" = '  +0.01%  '
    "E9" = '  +9.85%  '
    "E10" = '  -0.32%  '
    "E11" = '  +1.30%  '
    "E12" = '  +2.44%  '
    "E13" = '  +1.32%  '
    "E14" = '  +2.82%  '
    "E15" = '  +3.37%  '
    "E16" = '  +2.29%  '
    "E17" = '  +1.49%  '
    "E18" = '  +3.54%  '
    "E19" = '  +0.50%  '
    "E20" = '  +1.05%  '
    "E21" = '  +1.86%  '
    "E22" = '  +1.43%  '
    "E23" = '  +1.49%  '
    "E24" = '  +1.89%  '
    "E25" = '  -0.50%  '
    "E26" = '  -4.40%  '
    "E27" = '  -2.13%  '
    "E28" = '  +2.86%  '
    "E29" = '  +1.77%  '
    "E30" = '  +14.33%  '
    "E31" = '  +3.03%  '
    "E32" = '  -0.76%  '
    "E33" = '  +1.72%  '
    "E34" = '  +0.44%  '
    "E35" = '  +6.01%  '
    "E37" = '  +0.14%  '
    "E38" = '  +7.41%  '
    "E39" = '  -0.87%  '
    "E40" = '  +0.31%  '
    "E41" = '  +1.18%  '
    "E42" = '  -2.21%  '
    "E43" = '  +3.33%  '
    "E44" = '  +2.25%  '
    "E45" = '  +6.69%  '
    "E46" = '  +4.38%  '
    "E47" = '  +1.05%  '
    "E48" = '  +2.27%  '
    "E49" = '  -0.53%  '
    "E50" = '  +2.63%  '
    "E51" = '  +1.24%  '
}
foreach ($addr in $volumeUpdates.Keys) {
    $ws.Range($addr).Value2 = $volumeUpdates[$addr]
}
